# Apply "updates to contact and methods" edit to geneData.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet 1: studyData ---
$ws1 = $wb.Worksheets.Item("studyData")

# Row 2: disorder column (X) changes from "Dyslexia" to "reading"
$ws1.Range("X2").Value = "reading"

# Row 3: new study record (studyId 3) replacing the duplicate of row 2
$ws1.Range("A3").Value = 3
# B3 must stay text (looks numeric) -> force text format before writing
$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("B3").Value = "4000000"
$ws1.Range("E3").Value = "African"
$ws1.Range("X3").Value = "reading"

# --- Sheet 2: SNP_entryData ---
$ws2 = $wb.Worksheets.Item("SNP_entryData")

# New header column AC1: study_id
$ws2.Range("AC1").Value = "study_id"

# Row 2 becomes what used to be in row 3 (real SNP entry data)
$ws2.Range("A2").Value = 2
$ws2.Range("B2").Value = "xas"
$ws2.Range("C2").Value = 3.4
$ws2.Range("D2").Value = "cx1"
$ws2.Range("E2").Value = "X"
$ws2.Range("H2").Value = 5
$ws2.Range("I2").Value = 5
$ws2.Range("J2").Value = 5
$ws2.Range("K2").Value = "white"
$ws2.Range("L2").Value = "A"
$ws2.Range("M2").Value = "C"
$ws2.Range("N2").Value = 23
$ws2.Range("O2").Value = 23
$ws2.Range("P2").Value = 23
$ws2.Range("Q2").Value = "high"
$ws2.Range("R2").Value = "high"
$ws2.Range("S2").Value = "high.com"
$ws2.Range("T2").Value = "KROK"
# U2 must stay text (looks numeric) -> force text format before writing
$ws2.Range("U2").NumberFormat = "@"
$ws2.Range("U2").Value = "3"
$ws2.Range("V2").Value = "high"
$ws2.Range("W2").Value = "high"
$ws2.Range("Y2").Value = 10
$ws2.Range("Z2").Value = "high"
$ws2.Range("AA2").Value = "high.com"

# Row 3 is updated to a new gene entry
$ws2.Range("A3").Value = 3
$ws2.Range("B3").Value = "mee"
$ws2.Range("AB3").Value = "learning"
